$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (old row 4) entirely; the dataset now only has
# two event rows (rows 2-3), so the dimension shrinks from A1:L4 to A1:L3.
$ws.Rows("4:4").Delete()

# Row 2: new event data (No Seat Belt - FERNANDO ORNELAS)
$ws.Range("A2").Value = "281474991109864-1750183614996"
$ws.Range("B2").Value = "No Seat Belt"
$ws.Range("C2").Value = "2025-06-17T12:06:54.996"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "281474991109864"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "138"
$ws.Range("E2").ClearFormats()
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "51833957"
$ws.Range("F2").ClearFormats()
$ws.Range("G2").Value = "FERNANDO ORNELAS"
$ws.Range("H2").Value = 20.666377409
$ws.Range("I2").Value = -103.40913671
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109864/1750183612496/A5MTrC4t3W-camera-video-segment-driver-1750183614996.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSAAHM6GR4%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T150005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKT%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIAY0bsakobmrVnisLM2JlBXyLdh3La5Ypsdj5MWLhXrAAiEAmNVcKcqL6BmeVGprgDBES3RWF72IzNmKbaNmBs%2FvrBwq5gMIjf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDAsef%2BNlHp82B0CdNyq6A9yfq1k5ijE%2BhQze0XbSw2Vu3g9FtISCKgQbLJzrzxyg2pdnKilJNDO2Lv1Pt5Ph7%2B5OPg1lQ8Rik03miPkCfGmy0mC%2B8QqPhbrPxpLY9PVWs8HVGVWpPTS%2F57sQK2k7uoOwT61kxe4po5lSmiEnva9BzY1957gsoufupsbysCmbQGrZ6u%2FHYJVUWvWQgiiyMg26WE4g3QCPktE3msdl5ZWDkXjReaj5BarZhFnPjPw36ea4yxv4h8NRb8L%2BOxwSN7ipfNb%2BlsIM9B6MeGF93NcQ2bPYw6712g4CtJ6PZCnnp615zMLxdBuJDyzPqi9Xc8UqO1cwsLZw3o%2BUirIZtPyDgaCMC7AGXM83IzNnC9h%2FWvsRCfxc6YMf4vum9zUB03O1R9YeF2wOfUpTAlKMpJZnAWVJdWkGYVazxs03yo857FS%2F2RqtRcnmgYmnkKByvcm6%2BC1FeixFApx%2FeAZOnEAcg3eSn55CXrXj2Tw1J9iLhF0zblmw0QW4mShvD%2Fh%2Fa5iik2B6mvTufELpuZWMM8brJY8r9xhVdhwzg7unt%2BNELvoDqAsb7k4wj9B0NSVYNcBaPdVmKtNgfMAw2NLKwgY6pQFs3XweOVJ0mDN5OCWGiDUEO3rCgOZ6%2Fx8P28Lc44O2lshoPzAfZXmIRGHf8j7lO6dTGcB6yehIPZ8TbxrV3cfTMtBqT04S3N8xw5g6gD0kl1cqLDKXDCqo4yNvMW583PUynghEYYiLmsYvMR2EmrfZYG47saMinQMQSbjOFcoARPwytOxuNzYeuCfqR18PbdZTj6BMAm9RZm2vlLSext%2FQYefA7Pk%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2018%20Jun%202025%2023%3A00%3A05%20GMT&X-Amz-Signature=7d161eb05c024a7890a113a1909415f70e4ff982b00c6991e79ffbe6b20f9a92'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1750183612496/TAdN8NokD3-camera-video-segment-1750183614996.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSAAHM6GR4%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T150005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKT%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIAY0bsakobmrVnisLM2JlBXyLdh3La5Ypsdj5MWLhXrAAiEAmNVcKcqL6BmeVGprgDBES3RWF72IzNmKbaNmBs%2FvrBwq5gMIjf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDAsef%2BNlHp82B0CdNyq6A9yfq1k5ijE%2BhQze0XbSw2Vu3g9FtISCKgQbLJzrzxyg2pdnKilJNDO2Lv1Pt5Ph7%2B5OPg1lQ8Rik03miPkCfGmy0mC%2B8QqPhbrPxpLY9PVWs8HVGVWpPTS%2F57sQK2k7uoOwT61kxe4po5lSmiEnva9BzY1957gsoufupsbysCmbQGrZ6u%2FHYJVUWvWQgiiyMg26WE4g3QCPktE3msdl5ZWDkXjReaj5BarZhFnPjPw36ea4yxv4h8NRb8L%2BOxwSN7ipfNb%2BlsIM9B6MeGF93NcQ2bPYw6712g4CtJ6PZCnnp615zMLxdBuJDyzPqi9Xc8UqO1cwsLZw3o%2BUirIZtPyDgaCMC7AGXM83IzNnC9h%2FWvsRCfxc6YMf4vum9zUB03O1R9YeF2wOfUpTAlKMpJZnAWVJdWkGYVazxs03yo857FS%2F2RqtRcnmgYmnkKByvcm6%2BC1FeixFApx%2FeAZOnEAcg3eSn55CXrXj2Tw1J9iLhF0zblmw0QW4mShvD%2Fh%2Fa5iik2B6mvTufELpuZWMM8brJY8r9xhVdhwzg7unt%2BNELvoDqAsb7k4wj9B0NSVYNcBaPdVmKtNgfMAw2NLKwgY6pQFs3XweOVJ0mDN5OCWGiDUEO3rCgOZ6%2Fx8P28Lc44O2lshoPzAfZXmIRGHf8j7lO6dTGcB6yehIPZ8TbxrV3cfTMtBqT04S3N8xw5g6gD0kl1cqLDKXDCqo4yNvMW583PUynghEYYiLmsYvMR2EmrfZYG47saMinQMQSbjOFcoARPwytOxuNzYeuCfqR18PbdZTj6BMAm9RZm2vlLSext%2FQYefA7Pk%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2018%20Jun%202025%2023%3A00%3A05%20GMT&X-Amz-Signature=023794fc256ab8d162ab9d250bbd91e56ae9d70cca7be0b8b07ed2ba3461fe1f'

# Row 3: new event data (No Seat Belt - MARCOS BARBOSA)
$ws.Range("A3").Value = "281474991152736-1750182964183"
$ws.Range("B3").Value = "No Seat Belt"
$ws.Range("C3").Value = "2025-06-17T11:56:04.183"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "281474991152736"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "147"
$ws.Range("E3").ClearFormats()
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "51834001"
$ws.Range("F3").ClearFormats()
$ws.Range("G3").Value = "MARCOS BARBOSA"
$ws.Range("H3").Value = 20.67690181
$ws.Range("I3").Value = -103.36661835
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1750182961683/iBjUYEdqdu-camera-video-segment-driver-1750182964183.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSAAHM6GR4%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T150005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKT%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIAY0bsakobmrVnisLM2JlBXyLdh3La5Ypsdj5MWLhXrAAiEAmNVcKcqL6BmeVGprgDBES3RWF72IzNmKbaNmBs%2FvrBwq5gMIjf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDAsef%2BNlHp82B0CdNyq6A9yfq1k5ijE%2BhQze0XbSw2Vu3g9FtISCKgQbLJzrzxyg2pdnKilJNDO2Lv1Pt5Ph7%2B5OPg1lQ8Rik03miPkCfGmy0mC%2B8QqPhbrPxpLY9PVWs8HVGVWpPTS%2F57sQK2k7uoOwT61kxe4po5lSmiEnva9BzY1957gsoufupsbysCmbQGrZ6u%2FHYJVUWvWQgiiyMg26WE4g3QCPktE3msdl5ZWDkXjReaj5BarZhFnPjPw36ea4yxv4h8NRb8L%2BOxwSN7ipfNb%2BlsIM9B6MeGF93NcQ2bPYw6712g4CtJ6PZCnnp615zMLxdBuJDyzPqi9Xc8UqO1cwsLZw3o%2BUirIZtPyDgaCMC7AGXM83IzNnC9h%2FWvsRCfxc6YMf4vum9zUB03O1R9YeF2wOfUpTAlKMpJZnAWVJdWkGYVazxs03yo857FS%2F2RqtRcnmgYmnkKByvcm6%2BC1FeixFApx%2FeAZOnEAcg3eSn55CXrXj2Tw1J9iLhF0zblmw0QW4mShvD%2Fh%2Fa5iik2B6mvTufELpuZWMM8brJY8r9xhVdhwzg7unt%2BNELvoDqAsb7k4wj9B0NSVYNcBaPdVmKtNgfMAw2NLKwgY6pQFs3XweOVJ0mDN5OCWGiDUEO3rCgOZ6%2Fx8P28Lc44O2lshoPzAfZXmIRGHf8j7lO6dTGcB6yehIPZ8TbxrV3cfTMtBqT04S3N8xw5g6gD0kl1cqLDKXDCqo4yNvMW583PUynghEYYiLmsYvMR2EmrfZYG47saMinQMQSbjOFcoARPwytOxuNzYeuCfqR18PbdZTj6BMAm9RZm2vlLSext%2FQYefA7Pk%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2018%20Jun%202025%2023%3A00%3A05%20GMT&X-Amz-Signature=5676eec691bc4fc7e001901f6adfc7cfca44e974cfab62007a8c7b5c8b1bf6d1'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1750182961683/Foph7ZSROm-camera-video-segment-1750182964183.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSAAHM6GR4%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T150005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKT%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIAY0bsakobmrVnisLM2JlBXyLdh3La5Ypsdj5MWLhXrAAiEAmNVcKcqL6BmeVGprgDBES3RWF72IzNmKbaNmBs%2FvrBwq5gMIjf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDAsef%2BNlHp82B0CdNyq6A9yfq1k5ijE%2BhQze0XbSw2Vu3g9FtISCKgQbLJzrzxyg2pdnKilJNDO2Lv1Pt5Ph7%2B5OPg1lQ8Rik03miPkCfGmy0mC%2B8QqPhbrPxpLY9PVWs8HVGVWpPTS%2F57sQK2k7uoOwT61kxe4po5lSmiEnva9BzY1957gsoufupsbysCmbQGrZ6u%2FHYJVUWvWQgiiyMg26WE4g3QCPktE3msdl5ZWDkXjReaj5BarZhFnPjPw36ea4yxv4h8NRb8L%2BOxwSN7ipfNb%2BlsIM9B6MeGF93NcQ2bPYw6712g4CtJ6PZCnnp615zMLxdBuJDyzPqi9Xc8UqO1cwsLZw3o%2BUirIZtPyDgaCMC7AGXM83IzNnC9h%2FWvsRCfxc6YMf4vum9zUB03O1R9YeF2wOfUpTAlKMpJZnAWVJdWkGYVazxs03yo857FS%2F2RqtRcnmgYmnkKByvcm6%2BC1FeixFApx%2FeAZOnEAcg3eSn55CXrXj2Tw1J9iLhF0zblmw0QW4mShvD%2Fh%2Fa5iik2B6mvTufELpuZWMM8brJY8r9xhVdhwzg7unt%2BNELvoDqAsb7k4wj9B0NSVYNcBaPdVmKtNgfMAw2NLKwgY6pQFs3XweOVJ0mDN5OCWGiDUEO3rCgOZ6%2Fx8P28Lc44O2lshoPzAfZXmIRGHf8j7lO6dTGcB6yehIPZ8TbxrV3cfTMtBqT04S3N8xw5g6gD0kl1cqLDKXDCqo4yNvMW583PUynghEYYiLmsYvMR2EmrfZYG47saMinQMQSbjOFcoARPwytOxuNzYeuCfqR18PbdZTj6BMAm9RZm2vlLSext%2FQYefA7Pk%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2018%20Jun%202025%2023%3A00%3A05%20GMT&X-Amz-Signature=764becdc8856b6a59d84905bfb35579d3b05cad5050b5b4ac3d640c4a8a48e05'
